$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# --- Title-case Spanish prepositions in state/municipality names ---
$ws.Range("B4").Value = 'Rincón De Romos'
$ws.Range("B22").Value = 'Amatenango De La Frontera'
$ws.Range("B23").Value = 'Amatenango Del Valle'
$ws.Range("B26").Value = 'Bejucal De Ocampo'
$ws.Range("B31").Value = 'Chiapa De Corzo'
$ws.Range("B34").Value = 'Comitán De Domínguez'
$ws.Range("B52").Value = 'Mazapa De Madero'
$ws.Range("B56").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B62").Value = 'San Cristóbal De Las Casas'
$ws.Range("B63").Value = 'Santiago El Pinar'
$ws.Range("B103").Value = 'San Juan De Sabinas'
$ws.Range("B112").Value = 'Villa De Álvarez'
$ws.Range("A114").Value = 'Ciudad De México'
$ws.Range("B118").Value = 'Cuajimalpa De Morelos'
$ws.Range("A138").Value = 'Estado De México'
$ws.Range("B138").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B141").Value = 'Almoloya De Alquisiras'
$ws.Range("B142").Value = 'Almoloya De Juárez'
$ws.Range("B143").Value = 'Almoloya Del Río'
$ws.Range("B149").Value = 'Atizapán De Zaragoza'
$ws.Range("B154").Value = 'Chapa De Mota'
$ws.Range("B158").Value = 'Coacalco De Berriozábal'
$ws.Range("B163").Value = 'Ecatepec De Morelos'
$ws.Range("B170").Value = 'Ixtapan De La Sal'
$ws.Range("B179").Value = 'Naucalpan De Juárez'
$ws.Range("B186").Value = 'San Felipe Del Progreso'
$ws.Range("B187").Value = 'San Martín De Las Pirámides'
$ws.Range("B189").Value = 'Soyaniquilpan De Juárez'
$ws.Range("B198").Value = 'Tenango Del Aire'
$ws.Range("B199").Value = 'Tenango Del Valle'
$ws.Range("B208").Value = 'Tlalnepantla De Baz'
$ws.Range("B213").Value = 'Valle De Bravo'
$ws.Range("B214").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B215").Value = 'Villa De Allende'
$ws.Range("B216").Value = 'Villa Del Carbón'
$ws.Range("B226").Value = 'Apaseo El Alto'
$ws.Range("B227").Value = 'Apaseo El Grande'
$ws.Range("B233").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B240").Value = 'Purísima Del Rincón'
$ws.Range("B245").Value = 'San Luis De La Paz'
$ws.Range("B246").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B248").Value = 'Silao De La Victoria'
$ws.Range("B253").Value = 'Valle De Santiago'
$ws.Range("B257").Value = 'Acapulco De Juárez'
$ws.Range("B260").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B261").Value = 'Alcozauca De Guerrero'
$ws.Range("B265").Value = 'Atenango Del Río'
$ws.Range("B266").Value = 'Atlamajalcingo Del Monte'
$ws.Range("B268").Value = 'Atoyac De Álvarez'
$ws.Range("B269").Value = 'Ayutla De Los Libres'
$ws.Range("B272").Value = 'Buenavista De Cuéllar'
$ws.Range("B273").Value = 'Chilapa De Álvarez'
$ws.Range("B274").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B278").Value = 'Coyuca De Benítez'
$ws.Range("B279").Value = 'Coyuca De Catalán'
$ws.Range("B283").Value = 'Cuetzala Del Progreso'
$ws.Range("B284").Value = 'Cutzamala De Pinzón'
$ws.Range("B290").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B291").Value = 'Iguala De La Independencia'
$ws.Range("B293").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B294").Value = 'Zihuatanejo De Azueta'
$ws.Range("B296").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B299").Value = 'Mártir De Cuilapan'
$ws.Range("B311").Value = 'Taxco De Alarcón'
$ws.Range("B313").Value = 'Técpan De Galeana'
$ws.Range("B315").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B317").Value = 'Tixtla De Guerrero'
$ws.Range("B321").Value = 'Tlalixtaquilla De Maldonado'
$ws.Range("B322").Value = 'Tlapa De Comonfort'
$ws.Range("B334").Value = 'Agua Blanca De Iturbide'
$ws.Range("B338").Value = 'Atotonilco El Grande'
$ws.Range("B342").Value = 'Cuautepec De Hinojosa'
$ws.Range("B345").Value = 'Huasca De Ocampo'
$ws.Range("B348").Value = 'Huejutla De Reyes'
$ws.Range("B351").Value = 'Jacala De Ledezma'
$ws.Range("B356").Value = 'Mineral Del Chico'
$ws.Range("B357").Value = 'Mineral Del Monte'
$ws.Range("B358").Value = 'Mixquiahuala De Juárez'
$ws.Range("B359").Value = 'Molango De Escamilla'
$ws.Range("B360").Value = 'Pachuca De Soto'
$ws.Range("B362").Value = 'Progreso De Obregón'
$ws.Range("B369").Value = 'Tenango De Doria'
$ws.Range("B371").Value = 'Tepehuacán De Guerrero'
$ws.Range("B372").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B374").Value = 'Tezontepec De Aldama'
$ws.Range("B382").Value = 'Tula De Allende'
$ws.Range("B383").Value = 'Tulancingo De Bravo'
$ws.Range("B384").Value = 'Villa De Tezontepec'
$ws.Range("B387").Value = 'Zacualtipán De Ángeles'
$ws.Range("B390").Value = 'Acatlán De Juárez'
$ws.Range("B393").Value = 'Atotonilco El Alto'
$ws.Range("B395").Value = 'Autlán De Navarro'
$ws.Range("B403").Value = 'Encarnación De Díaz'
$ws.Range("B407").Value = 'Lagos De Moreno'
$ws.Range("B414").Value = 'San Juan De Los Lagos'
$ws.Range("B416").Value = 'San Martín De Bolaños'
$ws.Range("B417").Value = 'San Miguel El Alto'
$ws.Range("B420").Value = 'Tamazula De Gordiano'
$ws.Range("B423").Value = 'Teocuitatlán De Corona'
$ws.Range("B430").Value = 'Zapotlán El Grande'
$ws.Range("B449").Value = 'Coalcomán De Vázquez Pallares'
$ws.Range("B510").Value = 'Coatlán Del Río'
$ws.Range("B518").Value = 'Jonacatepec De Leandro Valle'
$ws.Range("B522").Value = 'Puente De Ixtla'
$ws.Range("B528").Value = 'Tetela Del Volcán'
$ws.Range("B530").Value = 'Tlaltizapán De Zapata'
$ws.Range("B538").Value = 'Zacualpan De Amilpas'
$ws.Range("B541").Value = 'Bahía De Banderas'
$ws.Range("B556").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B561").Value = 'Ayoquezco De Aldama'
$ws.Range("B564").Value = 'Capulálpam De Méndez'
$ws.Range("B566").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B567").Value = 'Ciénega De Zimatlán'
$ws.Range("B570").Value = 'Coicoyán De Las Flores'
$ws.Range("B573").Value = 'Constancia Del Rosario'
$ws.Range("B576").Value = 'Cuilápam De Guerrero'
$ws.Range("B577").Value = 'Cuyamecalco Villa De Zaragoza'
$ws.Range("B578").Value = 'El Barrio De La Soledad'
$ws.Range("B579").Value = 'Fresnillo De Trujano'
$ws.Range("B580").Value = 'Guadalupe De Ramírez'
$ws.Range("B582").Value = 'Guelatao De Juárez'
$ws.Range("B583").Value = 'Guevea De Humboldt'
$ws.Range("B584").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B585").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B586").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B587").Value = 'Huautla De Jiménez'
$ws.Range("B589").Value = 'Ixtlán De Juárez'
$ws.Range("B590").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B600").Value = 'Mariscala De Juárez'
$ws.Range("B601").Value = 'Mártires De Tacubaya'
$ws.Range("B604").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B605").Value = 'Mixistlán De La Reforma'
$ws.Range("B607").Value = 'Nejapa De Madero'
$ws.Range("B609").Value = 'Oaxaca De Juárez'
$ws.Range("B610").Value = 'Ocotlán De Morelos'
$ws.Range("B611").Value = 'Pinotepa De Don Luis'
$ws.Range("B613").Value = 'Putla Villa De Guerrero'
$ws.Range("B614").Value = 'Reforma De Pineda'
$ws.Range("B616").Value = 'Rojas De Cuauhtémoc'
$ws.Range("B632").Value = 'San Antonino El Alto'
$ws.Range("B634").Value = 'San Antonio De La Cal'
$ws.Range("B638").Value = 'San Baltazar Yatzachi El Bajo'
$ws.Range("B647").Value = 'San Felipe Jalapa De Díaz'
$ws.Range("B651").Value = 'San Francisco Del Mar'
$ws.Range("B672").Value = 'San José Del Progreso'
$ws.Range("B681").Value = 'San Juan Bautista Lo De Soto'
$ws.Range("B691").Value = 'San Juan Del Estado'
$ws.Range("B692").Value = 'San Juan Del Río'
$ws.Range("B749").Value = 'San Miguel Del Puerto'
$ws.Range("B751").Value = 'San Miguel El Grande'
$ws.Range("B772").Value = 'San Pablo Villa De Mitla'
$ws.Range("B777").Value = 'San Pedro El Alto'
$ws.Range("B795").Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Range("B796").Value = 'San Pedro Y San Pablo Tequixtepec'
$ws.Range("B821").Value = 'Santa Cruz De Bravo'
$ws.Range("B825").Value = 'Santa Cruz Tacache De Mina'
$ws.Range("B830").Value = 'Santa Inés Del Monte'
$ws.Range("B899").Value = 'Santo Domingo De Morelos'
$ws.Range("B918").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B919").Value = 'Tataltepec De Valdés'
$ws.Range("B920").Value = 'Teococuilco De Marcos Pérez'
$ws.Range("B921").Value = 'Teotitlán De Flores Magón'
$ws.Range("B922").Value = 'Teotitlán Del Valle'
$ws.Range("B924").Value = 'Tezoatlán De Segura Y Luna'
$ws.Range("B925").Value = 'Tlacolula De Matamoros'
$ws.Range("B926").Value = 'Tlalixtac De Cabrera'
$ws.Range("B927").Value = 'Totontepec Villa De Morelos'
$ws.Range("B929").Value = 'Villa De Etla'
$ws.Range("B930").Value = 'Villa De Tamazulápam Del Progreso'
$ws.Range("B931").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B932").Value = 'Villa De Zaachila'
$ws.Range("B934").Value = 'Villa Sola De Vega'
$ws.Range("B935").Value = 'Villa Talea De Castro'
$ws.Range("B938").Value = 'Zimatlán De Álvarez'
$ws.Range("B963").Value = 'Ayotoxco De Guerrero'
$ws.Range("B968").Value = 'Chalchicomula De Sesma'
$ws.Range("B978").Value = 'Chila De La Sal'
$ws.Range("B989").Value = 'Cuapiaxtla De Madero'
$ws.Range("B992").Value = 'Cuayuca De Andrade'
$ws.Range("B993").Value = 'Cuetzalan Del Progreso'
$ws.Range("B1009").Value = 'Huehuetlán El Chico'
$ws.Range("B1010").Value = 'Huehuetlán El Grande'
$ws.Range("B1015").Value = 'Huitzilan De Serdán'
$ws.Range("B1017").Value = 'Ixcamilpa De Guerrero'
$ws.Range("B1020").Value = 'Izúcar De Matamoros'
$ws.Range("B1029").Value = 'Los Reyes De Juárez'
$ws.Range("B1030").Value = 'Mazapiltepec De Juárez'
$ws.Range("B1042").Value = 'Palmar De Bravo'
$ws.Range("B1052").Value = 'San Diego La Mesa Tochimiltzingo'
$ws.Range("B1067").Value = 'San Nicolás De Los Ranchos'
$ws.Range("B1071").Value = 'San Salvador El Seco'
$ws.Range("B1072").Value = 'San Salvador El Verde'
$ws.Range("B1079").Value = 'Tecali De Herrera'
$ws.Range("B1087").Value = 'Tepanco De López'
$ws.Range("B1088").Value = 'Tepango De Rodríguez'
$ws.Range("B1089").Value = 'Tepatlaxco De Hidalgo'
$ws.Range("B1095").Value = 'Tepexi De Rodríguez'
$ws.Range("B1097").Value = 'Tepeyahualco De Cuauhtémoc'
$ws.Range("B1098").Value = 'Tetela De Ocampo'
$ws.Range("B1099").Value = 'Teteles De Avila Castillo'
$ws.Range("B1104").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B1116").Value = 'Totoltepec De Guerrero'
$ws.Range("B1121").Value = 'Xayacatlán De Bravo'
$ws.Range("B1127").Value = 'Xochitlán De Vicente Suárez'
$ws.Range("B1135").Value = 'Zapotitlán De Méndez'
$ws.Range("B1142").Value = 'Cadereyta De Montes'
$ws.Range("B1146").Value = 'Jalpan De Serra'
$ws.Range("B1147").Value = 'Landa De Matamoros'
$ws.Range("B1149").Value = 'Pinal De Amoles'
$ws.Range("B1151").Value = 'San Juan Del Río'
$ws.Range("B1159").Value = 'Axtla De Terrazas'
$ws.Range("B1160").Value = 'Ciudad Del Maíz'
$ws.Range("B1176").Value = 'Tanquián De Escobedo'
$ws.Range("B1178").Value = 'Villa De Reyes'
$ws.Range("B1223").Value = 'Soto La Marina'
$ws.Range("B1229").Value = 'Acuamanala De Miguel Hidalgo'
$ws.Range("B1236").Value = 'Contla De Juan Cuamatzi'
$ws.Range("B1243").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B1246").Value = 'Mazatecochco De José María Morelos'
$ws.Range("B1247").Value = 'Nanacamilpa De Mariano Arista'
$ws.Range("B1250").Value = 'Papalotla De Xicohténcatl'
$ws.Range("B1253").Value = 'San Pablo Del Monte'
$ws.Range("B1254").Value = 'Sanctórum De Lázaro Cárdenas'
$ws.Range("B1260").Value = 'Tepetitla De Lardizábal'
$ws.Range("B1263").Value = 'Tetla De La Solidaridad'
$ws.Range("B1275").Value = 'Ziltlaltépec De Trinidad Sánchez Santos'
$ws.Range("B1281").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B1284").Value = 'Amatlán De Los Reyes'
$ws.Range("B1292").Value = 'Boca Del Río'
$ws.Range("B1297").Value = 'Castillo De Teayo'
$ws.Range("B1299").Value = 'Cazones De Herrera'
$ws.Range("B1316").Value = 'Cosamaloapan De Carpio'
$ws.Range("B1330").Value = 'Hueyapan De Ocampo'
$ws.Range("B1331").Value = 'Ignacio De La Llave'
$ws.Range("B1335").Value = 'Ixhuatlán De Madero'
$ws.Range("B1336").Value = 'Ixhuatlán Del Sureste'
$ws.Range("B1345").Value = 'Juchique De Ferrer'
$ws.Range("B1348").Value = 'Landero Y Coss'
$ws.Range("B1350").Value = 'Lerdo De Tejada'
$ws.Range("B1355").Value = 'Martínez De La Torre'
$ws.Range("B1361").Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range("B1370").Value = 'Paso De Ovejas'
$ws.Range("B1371").Value = 'Paso Del Macho'
$ws.Range("B1375").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1382").Value = 'Sayula De Alemán'
$ws.Range("B1406").Value = 'Vega De Alatorre'
$ws.Range("B1415").Value = 'Zontecomatlán De López Y Fuentes'
$ws.Range("B1416").Value = 'Zozocolco De Hidalgo'

# --- Floating point literal adjustments (1-ULP recalculation bumps) ---
$ws.Range("D52").Value = 0.0009187892399571232
$ws.Range("D313").Value = 0.0009187892399571232
$ws.Range("D754").Value = 0.0009187892399571232
$ws.Range("D1009").Value = 0.0009187892399571232
$ws.Range("D1130").Value = 0.0009187892399571232
$ws.Range("D1242").Value = 0.0009187892399571232
$ws.Range("D1261").Value = 0.0009187892399571232
$ws.Range("D1282").Value = 0.0009187892399571232
$ws.Range("D868").Value = 0.009800418559542648

# --- Remove trailing footer/metadata rows 1439:1443 ---
$ws.Rows("1439:1443").Delete()
